# LOB1263.docx edit script
#
# The commit rotates the text content among several paragraphs/runs:
#   P6  (PT objectives)          <-  PT "Programa resumido" summary text
#   P7  (EN objectives, italic)  <-  EN "Programa resumido" summary text (italic)
#   P10 (teacher line)           <-  PT objectives text
#   P11 (EN summary, italic)     <-  EN objectives text (italic)
#   P13 (PT "Programa" summary)  <-  PT numbered programme list (9 items)
#   P16 "Método:" run            <-  "Critério:" run text (Aulas expositivas...)
#   P16 "Critério:" run          <-  "Norma de recuperação:" run text (NF=...)
#   P16 "Norma de recuperação:" run <- "Bibliografia" paragraph text (Estará apto...)
#   P18 (Bibliografia paragraph) <-  teacher line text (5840820 - Gustavo...)
#
# Paragraph indices stay stable across these edits (no paragraphs are
# added or removed), so we address the simple, single-run paragraphs by
# their fixed Paragraphs collection index.  The compound "Avaliação"
# paragraph (index 16) holds three bold labels ("Método: ", "Critério: ",
# "Norma de recuperação: ") each followed by a plain run; those plain
# runs are addressed with Find-derived character offsets, edited from
# right to left so earlier offsets stay valid.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: locate the character offset immediately AFTER a literal label
# (including its trailing ":"), then skip one extra character for the
# space that follows it (xml:space="preserve">Label: </w:t>), returning
# the start offset of the content run that follows the bold label.
# ---------------------------------------------------------------------
function Get-ContentStartAfterLabel([string]$label) {
    $r = $d.Content
    $ok = $r.Find.Execute($label, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find label: $label"
    }
    return $r.End + 1
}

# ---------------------------------------------------------------------
# Phase 1: simple whole-paragraph text swaps (stable Paragraphs index).
# ---------------------------------------------------------------------

$textE = "Eco inovação. Métricas da eco-inovação. Introdução ao Ciclo de vida do produto. Eco inovação na indústria. Estudo de casos de projetos de eco-inovação no Brasil. Métodos e ferramentas suporte do processo de eco-inovação. Identificação antecipada de falha como suporte a eco-inovação. TRIZ como resposta a eco-inovação. Proposta metodológica para soluções eco-inovadoras."
$textD = "Eco-innovation. Eco-innovation metrics. Introduction to products life-cycle. Eco-innovation in the industry.  Case study of Eco-innovation projects in Brazil. Methods and tools to support the process of Eco-innovation. Early identification of failure as support to Eco-innovation. TRIZ as a response to Eco-innovation. Methodological proposal for Eco-innovative solutions."
$textA = "A disciplina visa apresentar aos estudantes o conceito, tipos, modelos e sistemas de eco-inovação para o desenvolvimento da capacidade analítica e propositiva como competências profissionais nas áreas de inovação e sustentabilidade."
$textB = "The course aims to present to students the concept, types, models and systems of eco-innovation for the development of analytical and propositive capacity as professional competences in the areas of innovation and sustainability."
$textC = "5840820 - Gustavo Aristides Santana Martinez"
$textI = "Estará apto a efetuar a prova de reavaliação o aluno que tiver como média final na disciplina uma nota igual ou superior a três (3,0) e inferior a cinco (5,0), e tiver, no mínimo, 70% (setenta por cento) de frequência às aulas. O cálculo de uma média aritmética simples será feito com a nota da prova de reavaliação e a média final obtida pelo aluno na disciplina. Se esta média resultar em nota igual ou superior a cinco (5,0), o aluno será aprovado."

$br = [char]11
$textG = (
    "1. Eco-inovação: conceitos, fatores determinantes, barreiras, tipos de agentes eco-inovadores, categorias de eco inovações." + $br +
    "2. Métricas da eco-inovação: métricas de Andersen, métricas de Arundel & Kemp, métricas da OECD." + $br +
    "3. Introdução ao Ciclo de vida do produto: perspectiva analítica, análise da cadeia de produção, práticas de Green Supply Chain Management." + $br +
    "4. Eco inovação na indústria: química, agro alimentos, metal mecânica." + $br +
    "5. Estudo de casos de projetos de eco inovação no Brasil." + $br +
    "6. Métodos e ferramentas suporte do processo de eco-inovação: Eco-Compass, Eco-Ideation Tool, Value Mapping Tool, Design for Environment (DfE), EcoASIT, outros." + $br +
    "7. Identificação antecipada de falha como suporte a eco-inovação: o problema, o cenário, os recursos." + $br +
    "8. TRIZ como resposta a eco inovação: princípios inventivos, parâmetros de engenharia, matriz das contradições." + $br +
    "9. Proposta metodológica para soluções eco inovadoras na categoria tecnologias: definir, medir, analisar, criar"
)

$d.Paragraphs.Item(6).Range.Text = $textE
$d.Paragraphs.Item(7).Range.Text = $textD
$d.Paragraphs.Item(10).Range.Text = $textA
$d.Paragraphs.Item(11).Range.Text = $textB
$d.Paragraphs.Item(13).Range.Text = $textG
$d.Paragraphs.Item(18).Range.Text = $textC

# ---------------------------------------------------------------------
# Phase 2: the three plain runs inside the "Avaliação" paragraph.
# Edited right-to-left (Norma -> Critério -> Método) so the offsets for
# runs further to the left remain valid while we still need them.
# ---------------------------------------------------------------------

$textF = "Aulas expositivas, discussão de casos em sala de aula, painéis, debates, seminários, análise de vídeos e palestrantes externos."
$textH = "NF= (N1 + N2)/2" + $br + "Onde: NF = nota final; N = nota"

# "Norma de recuperação: " run -> becomes the "Estará apto..." text (textI)
$normaStart = Get-ContentStartAfterLabel("Norma de recuperação:")
$normaLabelEnd = $normaStart - 1
$afterNorma = $d.Paragraphs.Item(16).Range.End - 1   # paragraph mark excluded
$rNorma = $d.Range($normaStart, $afterNorma)
$rNorma.Text = $textI

# "Critério: " run -> becomes the NF=... text (textH). Its own trailing
# line break (<w:br/>) stays put, so stop just before it.
$criterioStart = Get-ContentStartAfterLabel("Critério:")
$rCriterioProbe = $d.Range($criterioStart, $normaLabelEnd)
$criterioText = $rCriterioProbe.Text
$criterioEnd = $criterioStart + $criterioText.Length - 1   # drop trailing <w:br/> char
$rCriterio = $d.Range($criterioStart, $criterioEnd)
$rCriterio.Text = $textH

# "Método: " run -> becomes the "Aulas expositivas..." text (textF). Its
# own trailing line break (<w:br/>) stays put, so stop just before it.
$metodoStart = Get-ContentStartAfterLabel("Método:")
$rMetodoProbe = $d.Range($metodoStart, $criterioStart - 1)
$metodoText = $rMetodoProbe.Text
$metodoEnd = $metodoStart + $metodoText.Length - 1   # drop trailing <w:br/> char
$rMetodo = $d.Range($metodoStart, $metodoEnd)
$rMetodo.Text = $textF

Write-Output "done"
